$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "date" column (C) for rows 6-23 as plain text strings,
# matching the target values from the diff.
$dates = @{
    6  = "2023-11-02"
    7  = "2023-11-03"
    8  = "2023-11-04"
    9  = "2023-11-04"
    10 = "2023-11-04"
    11 = "2023-11-05"
    12 = "2023-11-06"
    13 = "2023-11-07"
    14 = "2023-11-07"
    15 = "2023-11-07"
    16 = "2023-11-08"
    17 = "2023-11-08"
    18 = "2023-11-08"
    19 = "2023-11-08"
    20 = "2023-11-09"
    21 = "2023-11-10"
    22 = "2023-11-10"
    23 = "2023-11-10"
}

foreach ($row in $dates.Keys) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.NumberFormat = "@"
    $cell.Value = $dates[$row]
}
